$d = $word.ActiveDocument

# Locate the "Platform impact" bullet in the KEY ACHIEVEMENTS AND IMPACT
# section -- the new bullets are inserted immediately after it and before
# the following "TECHNICAL SKILLS" heading.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations*") {
        $target = $p
    }
}

if ($target -eq $null) {
    Write-Output "ERROR: anchor paragraph not found"
} else {
    $bullets = @(
        "• Real-time collaboration at national scale",
        "• Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%",
        "• Increased voter turnout prediction accuracy from 71% to 87%",
        "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
    )

    $current = $target
    foreach ($bulletText in $bullets) {
        $current.Range.InsertParagraphAfter()
        $current = $current.Next()
        $current.Range.Text = $bulletText
    }

    Write-Output "Inserted $($bullets.Count) achievement bullets after 'Platform impact' paragraph."
}
